$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Update existing odds values on row 2
# ---------------------------------------------------------------------------
$row2Updates = @{
    "J2" = 1.04
    "K2" = 12
    "L2" = 1.22
    "N2" = 1.73
    "R2" = 1.83
    "S2" = 1.83
    "W2" = 10
    "Z2" = 12
    "AB2" = 17
    "AD2" = 251
}
foreach ($addr in $row2Updates.Keys) {
    $ws.Range($addr).Value = $row2Updates[$addr]
}

# ---------------------------------------------------------------------------
# 2) Update existing odds values on row 3
# ---------------------------------------------------------------------------
$row3Updates = @{
    "G3" = 1.36
    "H3" = 5.25
    "I3" = 7
    "R3" = 1.7
    "S3" = 2.05
    "T3" = 10
    "W3" = 10
    "AB3" = 17
    "AC3" = 41
    "AD3" = 151
    "AI3" = 41
}
foreach ($addr in $row3Updates.Keys) {
    $ws.Range($addr).Value = $row3Updates[$addr]
}

# ---------------------------------------------------------------------------
# 3) Append new match rows (4-7)
# ---------------------------------------------------------------------------

# Row 4: Kuressaare - Tammeka (no odds available)
$ws.Range("A4").Value = "tW1KuYM5"
$ws.Range("B4").Value = "17/06/2025"
$ws.Range("C4").Value = "12:00"
$ws.Range("D4").Value = "ESTONIA - MEISTRILIIGA"
$ws.Range("E4").Value = "Kuressaare"
$ws.Range("F4").Value = "Tammeka"

# Row 5: Metta - BFC Daugavpils
$ws.Range("A5").Value = "dhDbPSze"
$ws.Range("B5").Value = "17/06/2025"
$ws.Range("C5").Value = "12:00"
$ws.Range("D5").Value = "LATVIA - VIRSLIGA"
$ws.Range("E5").Value = "Metta"
$ws.Range("F5").Value = "BFC Daugavpils"
$row5Odds = @{
    "G" = 2.95; "H" = 3.3; "I" = 2.12
    "N" = 1.93; "O" = 1.7; "P" = 1.37; "Q" = 2.5
    "T" = 7.5; "U" = 12; "V" = 9.25; "W" = 28; "X" = 21; "Y" = 29
    "Z" = 9; "AA" = 5.6; "AB" = 12.5; "AC" = 55; "AD" = 400
    "AE" = 6.2; "AF" = 8.25; "AG" = 7.6; "AH" = 16; "AI" = 14.5; "AJ" = 24
}
foreach ($col in $row5Odds.Keys) {
    $ws.Range($col + "5").Value = $row5Odds[$col]
}

# Row 6: Tukums 2000 - Jelgava
$ws.Range("A6").Value = "hK67N657"
$ws.Range("B6").Value = "17/06/2025"
$ws.Range("C6").Value = "13:00"
$ws.Range("D6").Value = "LATVIA - VIRSLIGA"
$ws.Range("E6").Value = "Tukums 2000"
$ws.Range("F6").Value = "Jelgava"
$row6Odds = @{
    "G" = 2.85; "H" = 3.35; "I" = 2.18
    "N" = 1.8; "O" = 1.8; "P" = 1.38; "Q" = 2.47
    "T" = 8; "U" = 12.5; "V" = 8.75; "W" = 27; "X" = 19.5; "Y" = 25
    "Z" = 10.25; "AA" = 5.7; "AB" = 11.5; "AC" = 45; "AD" = 300
    "AE" = 7; "AF" = 9.25; "AG" = 7.6; "AH" = 17; "AI" = 14; "AJ" = 21
}
foreach ($col in $row6Odds.Keys) {
    $ws.Range($col + "6").Value = $row6Odds[$col]
}

# Row 7: Super Nova - FK Liepaja
$ws.Range("A7").Value = "48DhoAcE"
$ws.Range("B7").Value = "17/06/2025"
$ws.Range("C7").Value = "14:00"
$ws.Range("D7").Value = "LATVIA - VIRSLIGA"
$ws.Range("E7").Value = "Super Nova"
$ws.Range("F7").Value = "FK Liepaja"
$row7Odds = @{
    "G" = 2.42; "H" = 3.25; "I" = 2.57
    "N" = 1.93; "O" = 1.7; "P" = 1.38; "Q" = 2.47
    "T" = 6.7; "U" = 9.75; "V" = 8; "W" = 20; "X" = 16.5; "Y" = 25
    "Z" = 9; "AA" = 5.5; "AB" = 12; "AC" = 55; "AD" = 350
    "AE" = 6.9; "AF" = 10.25; "AG" = 8.5; "AH" = 22; "AI" = 18; "AJ" = 26
}
foreach ($col in $row7Odds.Keys) {
    $ws.Range($col + "7").Value = $row7Odds[$col]
}
